$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-style format from A464 to the new date cells A465:A491
$ws.Cells.Item(464,1).Copy() | Out-Null
$ws.Range("A465:A491").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(465,1).Value2 = 44539
$ws.Cells.Item(465,2).Value2 = 7
$ws.Cells.Item(465,3).Value2 = 48
$ws.Cells.Item(465,4).Value2 = 442.1925380009213
$ws.Cells.Item(466,1).Value2 = 44540
$ws.Cells.Item(466,2).Value2 = 4
$ws.Cells.Item(466,3).Value2 = 33
$ws.Cells.Item(466,4).Value2 = 304.0073698756333
$ws.Cells.Item(467,1).Value2 = 44541
$ws.Cells.Item(467,2).Value2 = 0
$ws.Cells.Item(467,3).Value2 = 33
$ws.Cells.Item(467,4).Value2 = 304.0073698756333
$ws.Cells.Item(468,1).Value2 = 44542
$ws.Cells.Item(468,2).Value2 = 7
$ws.Cells.Item(468,3).Value2 = 30
$ws.Cells.Item(468,4).Value2 = 276.3703362505758
$ws.Cells.Item(469,1).Value2 = 44543
$ws.Cells.Item(469,2).Value2 = 8
$ws.Cells.Item(469,3).Value2 = 33
$ws.Cells.Item(469,4).Value2 = 304.0073698756333
$ws.Cells.Item(470,1).Value2 = 44544
$ws.Cells.Item(470,2).Value2 = 5
$ws.Cells.Item(470,3).Value2 = 31
$ws.Cells.Item(470,4).Value2 = 285.5826807922617
$ws.Cells.Item(471,1).Value2 = 44545
$ws.Cells.Item(471,2).Value2 = 0
$ws.Cells.Item(471,3).Value2 = 31
$ws.Cells.Item(471,4).Value2 = 285.5826807922617
$ws.Cells.Item(472,1).Value2 = 44546
$ws.Cells.Item(472,2).Value2 = 3
$ws.Cells.Item(472,3).Value2 = 27
$ws.Cells.Item(472,4).Value2 = 248.7333026255182
$ws.Cells.Item(473,1).Value2 = 44547
$ws.Cells.Item(473,2).Value2 = 5
$ws.Cells.Item(473,3).Value2 = 28
$ws.Cells.Item(473,4).Value2 = 257.945647167204
$ws.Cells.Item(474,1).Value2 = 44548
$ws.Cells.Item(474,2).Value2 = 3
$ws.Cells.Item(474,3).Value2 = 31
$ws.Cells.Item(474,4).Value2 = 285.5826807922617
$ws.Cells.Item(475,1).Value2 = 44550
$ws.Cells.Item(475,2).Value2 = 7
$ws.Cells.Item(475,3).Value2 = 31
$ws.Cells.Item(475,4).Value2 = 285.5826807922617
$ws.Cells.Item(476,1).Value2 = 44551
$ws.Cells.Item(476,2).Value2 = 5
$ws.Cells.Item(476,3).Value2 = 28
$ws.Cells.Item(476,4).Value2 = 257.945647167204
$ws.Cells.Item(477,1).Value2 = 44552
$ws.Cells.Item(477,2).Value2 = 0
$ws.Cells.Item(477,3).Value2 = 23
$ws.Cells.Item(477,4).Value2 = 211.8839244587748
$ws.Cells.Item(478,1).Value2 = 44553
$ws.Cells.Item(478,2).Value2 = 7
$ws.Cells.Item(478,3).Value2 = 30
$ws.Cells.Item(478,4).Value2 = 276.3703362505758
$ws.Cells.Item(479,1).Value2 = 44554
$ws.Cells.Item(479,2).Value2 = 2
$ws.Cells.Item(479,3).Value2 = 29
$ws.Cells.Item(479,4).Value2 = 267.1579917088899
$ws.Cells.Item(480,1).Value2 = 44555
$ws.Cells.Item(480,2).Value2 = 4
$ws.Cells.Item(480,3).Value2 = 28
$ws.Cells.Item(480,4).Value2 = 257.945647167204
$ws.Cells.Item(481,1).Value2 = 44556
$ws.Cells.Item(481,2).Value2 = 5
$ws.Cells.Item(481,3).Value2 = 30
$ws.Cells.Item(481,4).Value2 = 276.3703362505758
$ws.Cells.Item(482,1).Value2 = 44557
$ws.Cells.Item(482,2).Value2 = 7
$ws.Cells.Item(482,3).Value2 = 30
$ws.Cells.Item(482,4).Value2 = 276.3703362505758
$ws.Cells.Item(483,1).Value2 = 44558
$ws.Cells.Item(483,2).Value2 = 20
$ws.Cells.Item(483,3).Value2 = 45
$ws.Cells.Item(483,4).Value2 = 414.5555043758637
$ws.Cells.Item(484,1).Value2 = 44559
$ws.Cells.Item(484,2).Value2 = 2
$ws.Cells.Item(484,3).Value2 = 47
$ws.Cells.Item(484,4).Value2 = 432.9801934592354
$ws.Cells.Item(485,1).Value2 = 44560
$ws.Cells.Item(485,2).Value2 = 8
$ws.Cells.Item(485,3).Value2 = 48
$ws.Cells.Item(485,4).Value2 = 442.1925380009213
$ws.Cells.Item(486,1).Value2 = 44561
$ws.Cells.Item(486,2).Value2 = 17
$ws.Cells.Item(486,3).Value2 = 63
$ws.Cells.Item(486,4).Value2 = 580.3777061262091
$ws.Cells.Item(487,1).Value2 = 44562
$ws.Cells.Item(487,2).Value2 = 12
$ws.Cells.Item(487,3).Value2 = 71
$ws.Cells.Item(487,4).Value2 = 654.076462459696
$ws.Cells.Item(488,1).Value2 = 44563
$ws.Cells.Item(488,2).Value2 = 13
$ws.Cells.Item(488,3).Value2 = 79
$ws.Cells.Item(488,4).Value2 = 727.7752187931828
$ws.Cells.Item(489,1).Value2 = 44564
$ws.Cells.Item(489,2).Value2 = 17
$ws.Cells.Item(489,3).Value2 = 89
$ws.Cells.Item(489,4).Value2 = 819.8986642100415
$ws.Cells.Item(490,1).Value2 = 44565
$ws.Cells.Item(490,2).Value2 = 6
$ws.Cells.Item(490,3).Value2 = 75
$ws.Cells.Item(490,4).Value2 = 690.9258406264395
$ws.Cells.Item(491,1).Value2 = 44566
$ws.Cells.Item(491,2).Value2 = 25
$ws.Cells.Item(491,3).Value2 = 98
$ws.Cells.Item(491,4).Value2 = 902.8097650852142

